$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRA_data")

$ws.Cells.Item(1, 18).Value = "srr_accession"
$ws.Cells.Item(8, 18).Value = "SRR24914067"
$ws.Cells.Item(9, 18).Value = "SRR24914092"
$ws.Cells.Item(11, 18).Value = "SRR25997770"
$ws.Cells.Item(12, 18).Value = "SRR24914063"
$ws.Cells.Item(13, 18).Value = "SRR24914055"
$ws.Cells.Item(15, 18).Value = "SRR24914066"
$ws.Cells.Item(17, 18).Value = "SRR24914091"
$ws.Cells.Item(19, 18).Value = "SRR24914065"
$ws.Cells.Item(21, 18).Value = "SRR24914090"
$ws.Cells.Item(23, 18).Value = "SRR24914064"
$ws.Cells.Item(25, 18).Value = "SRR24914089"
$ws.Cells.Item(27, 18).Value = "SRR24914062"
$ws.Cells.Item(29, 18).Value = "SRR25997769"
$ws.Cells.Item(30, 18).Value = "SRR24914054"
$ws.Cells.Item(31, 18).Value = "SRR24914060"
$ws.Cells.Item(33, 18).Value = "SRR24914053"
$ws.Cells.Item(35, 18).Value = "SRR24914050"
$ws.Cells.Item(37, 18).Value = "SRR24914057"
$ws.Cells.Item(39, 18).Value = "SRR24914021"
$ws.Cells.Item(41, 18).Value = "SRR24914048"
$ws.Cells.Item(43, 18).Value = "SRR24914023"
$ws.Cells.Item(44, 18).Value = "SRR24914043"
$ws.Cells.Item(46, 18).Value = "SRR24914056"
$ws.Cells.Item(47, 18).Value = "SRR25997762"
$ws.Cells.Item(48, 18).Value = "SRR24914046"
$ws.Cells.Item(50, 18).Value = "SRR24914116"
$ws.Cells.Item(52, 18).Value = "SRR25997761"
$ws.Cells.Item(53, 18).Value = "SRR24914038"
$ws.Cells.Item(54, 18).Value = "SRR24914114"
$ws.Cells.Item(56, 18).Value = "SRR24914036"
$ws.Cells.Item(58, 18).Value = "SRR24914111"
$ws.Cells.Item(60, 18).Value = "SRR25997760"
$ws.Cells.Item(61, 18).Value = "SRR24914032"
$ws.Cells.Item(62, 18).Value = "SRR25997759"
$ws.Cells.Item(63, 18).Value = "SRR24914108"
$ws.Cells.Item(64, 18).Value = "SRR24914030"
$ws.Cells.Item(66, 18).Value = "SRR24914103"
$ws.Cells.Item(68, 18).Value = "SRR25997758"
$ws.Cells.Item(69, 18).Value = "SRR24914025"
$ws.Cells.Item(70, 18).Value = "SRR24914105"
$ws.Cells.Item(72, 18).Value = "SRR24914027"
$ws.Cells.Item(74, 18).Value = "SRR24914106"
$ws.Cells.Item(76, 18).Value = "SRR24914028"
$ws.Cells.Item(78, 18).Value = "SRR24914107"
$ws.Cells.Item(80, 18).Value = "SRR24914029"
$ws.Cells.Item(82, 18).Value = "SRR24914102"
$ws.Cells.Item(84, 18).Value = "SRR24914024"
$ws.Cells.Item(86, 18).Value = "SRR24914097"
$ws.Cells.Item(88, 18).Value = "SRR24914015"
$ws.Cells.Item(90, 18).Value = "SRR24914095"
$ws.Cells.Item(92, 18).Value = "SRR24914013"
$ws.Cells.Item(94, 18).Value = "SRR24914059"
$ws.Cells.Item(96, 18).Value = "SRR24914052"
$ws.Cells.Item(104, 18).Value = "SRR24914123"
$ws.Cells.Item(105, 18).Value = "SRR24913982"
$ws.Cells.Item(108, 18).Value = "SRR24913985"
$ws.Cells.Item(109, 18).Value = "SRR24913978"
$ws.Cells.Item(111, 18).Value = "SRR24914122"
$ws.Cells.Item(113, 18).Value = "SRR24913981"
$ws.Cells.Item(115, 18).Value = "SRR24914007"
$ws.Cells.Item(117, 18).Value = "SRR24913980"
$ws.Cells.Item(119, 18).Value = "SRR24913996"
$ws.Cells.Item(121, 18).Value = "SRR24913979"
$ws.Cells.Item(123, 18).Value = "SRR24913974"
$ws.Cells.Item(126, 18).Value = "SRR24913977"
$ws.Cells.Item(127, 18).Value = "SRR24914080"
$ws.Cells.Item(129, 18).Value = "SRR24913976"
$ws.Cells.Item(131, 18).Value = "SRR24913972"
$ws.Cells.Item(133, 18).Value = "SRR24914061"
$ws.Cells.Item(135, 18).Value = "SRR24914110"
$ws.Cells.Item(137, 18).Value = "SRR24913970"
$ws.Cells.Item(139, 18).Value = "SRR24914034"
$ws.Cells.Item(140, 18).Value = "SRR24913966"
$ws.Cells.Item(142, 18).Value = "SRR24914088"
$ws.Cells.Item(144, 18).Value = "SRR24913968"
$ws.Cells.Item(146, 18).Value = "SRR24914006"
$ws.Cells.Item(149, 18).Value = "SRR24914083"
$ws.Cells.Item(150, 18).Value = "SRR24914004"
$ws.Cells.Item(152, 18).Value = "SRR24914085"
$ws.Cells.Item(154, 18).Value = "SRR24914001"
$ws.Cells.Item(157, 18).Value = "SRR24913956"
$ws.Cells.Item(159, 18).Value = "SRR24913999"
$ws.Cells.Item(160, 18).Value = "SRR24913958"
$ws.Cells.Item(162, 18).Value = "SRR24913993"
$ws.Cells.Item(165, 18).Value = "SRR24914070"
$ws.Cells.Item(166, 18).Value = "SRR24913995"
$ws.Cells.Item(168, 18).Value = "SRR24913962"
$ws.Cells.Item(170, 18).Value = "SRR24913997"
$ws.Cells.Item(172, 18).Value = "SRR24913961"
$ws.Cells.Item(174, 18).Value = "SRR24913998"
$ws.Cells.Item(176, 18).Value = "SRR24913960"
$ws.Cells.Item(178, 18).Value = "SRR24913992"
$ws.Cells.Item(180, 18).Value = "SRR24914071"
$ws.Cells.Item(182, 18).Value = "SRR24913988"
$ws.Cells.Item(184, 18).Value = "SRR24914075"
$ws.Cells.Item(186, 18).Value = "SRR24913986"
$ws.Cells.Item(188, 18).Value = "SRR24914079"
$ws.Cells.Item(190, 18).Value = "SRR24913959"
$ws.Cells.Item(192, 18).Value = "SRR24913975"
